$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: correct the D (Fecha) timestamp on the last existing run of rows
#         (rows 548-561) -- tiny floating point correction introduced by the
#         upstream "Actualizar" automation.
# ---------------------------------------------------------------------------
for ($r = 548; $r -le 561; $r++) {
    $ws.Cells.Item($r, 4).Value = 44232.49295149306
}

# ---------------------------------------------------------------------------
# Step 2: append a brand new 14-row "run" (rows 562-575) with the next
#         availability check, reusing the same Name/URL pattern that repeats
#         throughout the sheet.
# ---------------------------------------------------------------------------
$names = @("Odoo","Blackbox","PowerBI","Dropbox","Odoo","GEE","UtilidadesOdoo","Filtros Dashboard","MapStore","GeoServer","Tomcat","Shiny","Github","EZ Exporter")
$urls  = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)

$newDate  = 44232.51400987701
$startRow = 562

for ($i = 0; $i -lt 14; $i++) {
    $r = $startRow + $i

    $ws.Cells.Item($r, 1).Value = $names[$i]
    $ws.Cells.Item($r, 3).Value = "Disponible"
    $ws.Cells.Item($r, 4).Value = $newDate
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $h = $ws.Hyperlinks.Add($ws.Cells.Item($r, 2), $urls[$i])
    $ws.Cells.Item($r, 2).Value = $names[$i]
    $ws.Cells.Item($r, 2).Style = "Hyperlink"

    # The "MapStore" row links to a URL with a trailing "#/" fragment, which
    # OOXML encodes as a relationship target without the fragment plus a
    # location="/" attribute on the <hyperlink> element.
    if ($r -eq 570) {
        $h.SubAddress = "/"
    }
}
